$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.604.05'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.796.46'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.97'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.560'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.30%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.98'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.297'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.76%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.056.87'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.827.26'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.06'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.637'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.580.79'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.94'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '247.59'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0802'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.37'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.52%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.17'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '167.49'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.31'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.60'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.01'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.10'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +11.61%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0525'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.24'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.427.97'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.60'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +8.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.672'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.84%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0193'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.26%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.06'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.70%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.99%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.59%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.936'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.33%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.76'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.73'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0527'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.47%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.86%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.68%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '106.07'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -5.45%  '
